$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Ccl12"
$ws.Range("C2").Value2 = "Ccr5"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 1.005755333333333
$ws.Range("H2").Value2 = 3.017266
$ws.Range("I2").Value2 = 0.01048729000197281
$ws.Range("J2").Value2 = 0.01048729000197281
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.9442423333333334
$ws.Range("N2").Value2 = 2.832727
$ws.Range("O2").Value2 = 0.006848500623481535
$ws.Range("P2").Value2 = 0.006848500623481536
$ws.Range("Q2").Value2 = 0.9496767627091112
$ws.Range("R2").Value2 = 8.547090864382
$ws.Range("S2").Value2 = 0.00007182221211714243
$ws.Range("T2").Value2 = 0.00007182221211714243

$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Ccl12"
$ws.Range("C3").Value2 = "Ccr5"
$ws.Range("D3").Value2 = "Inflammatory-Mac"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 1.005755333333333
$ws.Range("H3").Value2 = 3.017266
$ws.Range("I3").Value2 = 0.01048729000197281
$ws.Range("J3").Value2 = 0.01048729000197281
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 82.477727
$ws.Range("N3").Value2 = 247.433181
$ws.Range("O3").Value2 = 0.5982031781913751
$ws.Range("P3").Value2 = 0.5982031781913751
$ws.Range("Q3").Value2 = 82.95241381146067
$ws.Range("R3").Value2 = 746.571724303146
$ws.Range("S3").Value2 = 0.006273530209794765
$ws.Range("T3").Value2 = 0.006273530209794765

$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Ccl12"
$ws.Range("C4").Value2 = "Ccr5"
$ws.Range("D4").Value2 = "Neutrophils"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 1.005755333333333
$ws.Range("H4").Value2 = 3.017266
$ws.Range("I4").Value2 = 0.01048729000197281
$ws.Range("J4").Value2 = 0.01048729000197281
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 9.766934000000001
$ws.Range("N4").Value2 = 29.300802
$ws.Range("O4").Value2 = 0.07083865150630789
$ws.Range("P4").Value2 = 0.07083865150630789
$ws.Range("Q4").Value2 = 9.823145960814667
$ws.Range("R4").Value2 = 88.40831364733199
$ws.Range("S4").Value2 = 0.0007429054816953385
$ws.Range("T4").Value2 = 0.0007429054816953385

$ws.Range("A5").Value2 = "ECs"
$ws.Range("B5").Value2 = "Ccl12"
$ws.Range("C5").Value2 = "Ccr5"
$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 1.005755333333333
$ws.Range("H5").Value2 = 3.017266
$ws.Range("I5").Value2 = 0.01048729000197281
$ws.Range("J5").Value2 = 0.01048729000197281
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 44.68687199999999
$ws.Range("N5").Value2 = 134.060616
$ws.Range("O5").Value2 = 0.3241096696788354
$ws.Range("P5").Value2 = 0.3241096696788355
$ws.Range("Q5").Value2 = 44.94405984398399
$ws.Range("R5").Value2 = 404.4965385958559
$ws.Range("S5").Value2 = 0.003399032098365559
$ws.Range("T5").Value2 = 0.00339903209836556

$ws.Range("A6").Value2 = "Inflammatory-Mac"
$ws.Range("B6").Value2 = "Ccl12"
$ws.Range("C6").Value2 = "Ccr5"
$ws.Range("D6").Value2 = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 38.755371
$ws.Range("H6").Value2 = 116.266113
$ws.Range("I6").Value2 = 0.4041130097356814
$ws.Range("J6").Value2 = 0.4041130097356814
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 0.9442423333333334
$ws.Range("N6").Value2 = 2.832727
$ws.Range("O6").Value2 = 0.006848500623481535
$ws.Range("P6").Value2 = 0.006848500623481536
$ws.Range("Q6").Value2 = 36.594461942239
$ws.Range("R6").Value2 = 329.350157480151
$ws.Range("S6").Value2 = 0.002767568199131814
$ws.Range("T6").Value2 = 0.002767568199131814

$ws.Range("A7").Value2 = "Inflammatory-Mac"
$ws.Range("B7").Value2 = "Ccl12"
$ws.Range("C7").Value2 = "Ccr5"
$ws.Range("D7").Value2 = "Inflammatory-Mac"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 38.755371
$ws.Range("H7").Value2 = 116.266113
$ws.Range("I7").Value2 = 0.4041130097356814
$ws.Range("J7").Value2 = 0.4041130097356814
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 82.477727
$ws.Range("N7").Value2 = 247.433181
$ws.Range("O7").Value2 = 0.5982031781913751
$ws.Range("P7").Value2 = 0.5982031781913751
$ws.Range("Q7").Value2 = 3196.454909121717
$ws.Range("R7").Value2 = 28768.09418209545
$ws.Range("S7").Value2 = 0.2417416867723667
$ws.Range("T7").Value2 = 0.2417416867723667

$ws.Range("A8").Value2 = "Inflammatory-Mac"
$ws.Range("B8").Value2 = "Ccl12"
$ws.Range("C8").Value2 = "Ccr5"
$ws.Range("D8").Value2 = "Neutrophils"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 38.755371
$ws.Range("H8").Value2 = 116.266113
$ws.Range("I8").Value2 = 0.4041130097356814
$ws.Range("J8").Value2 = 0.4041130097356814
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 9.766934000000001
$ws.Range("N8").Value2 = 29.300802
$ws.Range("O8").Value2 = 0.07083865150630789
$ws.Range("P8").Value2 = 0.07083865150630789
$ws.Range("Q8").Value2 = 378.521150702514
$ws.Range("R8").Value2 = 3406.690356322626
$ws.Range("S8").Value2 = 0.02862682066583114
$ws.Range("T8").Value2 = 0.02862682066583114

$ws.Range("A9").Value2 = "Inflammatory-Mac"
$ws.Range("B9").Value2 = "Ccl12"
$ws.Range("C9").Value2 = "Ccr5"
$ws.Range("D9").Value2 = "Resolving-Mac"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 38.755371
$ws.Range("H9").Value2 = 116.266113
$ws.Range("I9").Value2 = 0.4041130097356814
$ws.Range("J9").Value2 = 0.4041130097356814
$ws.Range("K9").Value2 = 3
$ws.Range("L9").Value2 = 1
$ws.Range("M9").Value2 = 44.68687199999999
$ws.Range("N9").Value2 = 134.060616
$ws.Range("O9").Value2 = 0.3241096696788354
$ws.Range("P9").Value2 = 0.3241096696788355
$ws.Range("Q9").Value2 = 1731.856303189512
$ws.Range("R9").Value2 = 15586.70672870561
$ws.Range("S9").Value2 = 0.1309769340983517
$ws.Range("T9").Value2 = 0.1309769340983517

$ws.Range("A10").Value2 = "Neutrophils"
$ws.Range("B10").Value2 = "Ccl12"
$ws.Range("C10").Value2 = "Ccr5"
$ws.Range("D10").Value2 = "ECs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 7.988471999999999
$ws.Range("H10").Value2 = 23.965416
$ws.Range("I10").Value2 = 0.08329801469605898
$ws.Range("J10").Value2 = 0.08329801469605898
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 0.9442423333333334
$ws.Range("N10").Value2 = 2.832727
$ws.Range("O10").Value2 = 0.006848500623481535
$ws.Range("P10").Value2 = 0.006848500623481536
$ws.Range("Q10").Value2 = 7.543053441047999
$ws.Range("R10").Value2 = 67.887480969432
$ws.Range("S10").Value2 = 0.000570466505580734
$ws.Range("T10").Value2 = 0.000570466505580734

$ws.Range("A11").Value2 = "Neutrophils"
$ws.Range("B11").Value2 = "Ccl12"
$ws.Range("C11").Value2 = "Ccr5"
$ws.Range("D11").Value2 = "Inflammatory-Mac"
$ws.Range("E11").Value2 = 3
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 7.988471999999999
$ws.Range("H11").Value2 = 23.965416
$ws.Range("I11").Value2 = 0.08329801469605898
$ws.Range("J11").Value2 = 0.08329801469605898
$ws.Range("K11").Value2 = 3
$ws.Range("L11").Value2 = 1
$ws.Range("M11").Value2 = 82.477727
$ws.Range("N11").Value2 = 247.433181
$ws.Range("O11").Value2 = 0.5982031781913751
$ws.Range("P11").Value2 = 0.5982031781913751
$ws.Range("Q11").Value2 = 658.8710127631439
$ws.Range("R11").Value2 = 5929.839114868295
$ws.Range("S11").Value2 = 0.04982913712821435
$ws.Range("T11").Value2 = 0.04982913712821435

$ws.Range("A12").Value2 = "Neutrophils"
$ws.Range("B12").Value2 = "Ccl12"
$ws.Range("C12").Value2 = "Ccr5"
$ws.Range("D12").Value2 = "Neutrophils"
$ws.Range("E12").Value2 = 3
$ws.Range("F12").Value2 = 1
$ws.Range("G12").Value2 = 7.988471999999999
$ws.Range("H12").Value2 = 23.965416
$ws.Range("I12").Value2 = 0.08329801469605898
$ws.Range("J12").Value2 = 0.08329801469605898
$ws.Range("K12").Value2 = 3
$ws.Range("L12").Value2 = 1
$ws.Range("M12").Value2 = 9.766934000000001
$ws.Range("N12").Value2 = 29.300802
$ws.Range("O12").Value2 = 0.07083865150630789
$ws.Range("P12").Value2 = 0.07083865150630789
$ws.Range("Q12").Value2 = 78.02287878484799
$ws.Range("R12").Value2 = 702.2059090636319
$ws.Range("S12").Value2 = 0.005900719034221435
$ws.Range("T12").Value2 = 0.005900719034221435

$ws.Range("A13").Value2 = "Neutrophils"
$ws.Range("B13").Value2 = "Ccl12"
$ws.Range("C13").Value2 = "Ccr5"
$ws.Range("D13").Value2 = "Resolving-Mac"
$ws.Range("E13").Value2 = 3
$ws.Range("F13").Value2 = 1
$ws.Range("G13").Value2 = 7.988471999999999
$ws.Range("H13").Value2 = 23.965416
$ws.Range("I13").Value2 = 0.08329801469605898
$ws.Range("J13").Value2 = 0.08329801469605898
$ws.Range("K13").Value2 = 3
$ws.Range("L13").Value2 = 1
$ws.Range("M13").Value2 = 44.68687199999999
$ws.Range("N13").Value2 = 134.060616
$ws.Range("O13").Value2 = 0.3241096696788354
$ws.Range("P13").Value2 = 0.3241096696788355
$ws.Range("Q13").Value2 = 356.9798257395839
$ws.Range("R13").Value2 = 3212.818431656255
$ws.Range("S13").Value2 = 0.02699769202804245
$ws.Range("T13").Value2 = 0.02699769202804246

$ws.Range("A14").Value2 = "Resolving-Mac"
$ws.Range("B14").Value2 = "Ccl12"
$ws.Range("C14").Value2 = "Ccr5"
$ws.Range("D14").Value2 = "ECs"
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 48.15271133333334
$ws.Range("H14").Value2 = 144.458134
$ws.Range("I14").Value2 = 0.5021016855662869
$ws.Range("J14").Value2 = 0.5021016855662868
$ws.Range("K14").Value2 = 3
$ws.Range("L14").Value2 = 1
$ws.Range("M14").Value2 = 0.9442423333333334
$ws.Range("N14").Value2 = 2.832727
$ws.Range("O14").Value2 = 0.006848500623481535
$ws.Range("P14").Value2 = 0.006848500623481536
$ws.Range("Q14").Value2 = 45.46782850571311
$ws.Range("R14").Value2 = 409.210456551418
$ws.Range("S14").Value2 = 0.003438643706651845
$ws.Range("T14").Value2 = 0.003438643706651845

$ws.Range("A15").Value2 = "Resolving-Mac"
$ws.Range("B15").Value2 = "Ccl12"
$ws.Range("C15").Value2 = "Ccr5"
$ws.Range("D15").Value2 = "Inflammatory-Mac"
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 48.15271133333334
$ws.Range("H15").Value2 = 144.458134
$ws.Range("I15").Value2 = 0.5021016855662869
$ws.Range("J15").Value2 = 0.5021016855662868
$ws.Range("K15").Value2 = 3
$ws.Range("L15").Value2 = 1
$ws.Range("M15").Value2 = 82.477727
$ws.Range("N15").Value2 = 247.433181
$ws.Range("O15").Value2 = 0.5982031781913751
$ws.Range("P15").Value2 = 0.5982031781913751
$ws.Range("Q15").Value2 = 3971.526179660473
$ws.Range("R15").Value2 = 35743.73561694425
$ws.Range("S15").Value2 = 0.3003588240809993
$ws.Range("T15").Value2 = 0.3003588240809992

$ws.Range("A16").Value2 = "Resolving-Mac"
$ws.Range("B16").Value2 = "Ccl12"
$ws.Range("C16").Value2 = "Ccr5"
$ws.Range("D16").Value2 = "Neutrophils"
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 48.15271133333334
$ws.Range("H16").Value2 = 144.458134
$ws.Range("I16").Value2 = 0.5021016855662869
$ws.Range("J16").Value2 = 0.5021016855662868
$ws.Range("K16").Value2 = 3
$ws.Range("L16").Value2 = 1
$ws.Range("M16").Value2 = 9.766934000000001
$ws.Range("N16").Value2 = 29.300802
$ws.Range("O16").Value2 = 0.07083865150630789
$ws.Range("P16").Value2 = 0.07083865150630789
$ws.Range("Q16").Value2 = 470.3043535137188
$ws.Range("R16").Value2 = 4232.739181623469
$ws.Range("S16").Value2 = 0.03556820632455997
$ws.Range("T16").Value2 = 0.03556820632455997

$ws.Range("A17").Value2 = "Resolving-Mac"
$ws.Range("B17").Value2 = "Ccl12"
$ws.Range("C17").Value2 = "Ccr5"
$ws.Range("D17").Value2 = "Resolving-Mac"
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 48.15271133333334
$ws.Range("H17").Value2 = 144.458134
$ws.Range("I17").Value2 = 0.5021016855662869
$ws.Range("J17").Value2 = 0.5021016855662868
$ws.Range("K17").Value2 = 3
$ws.Range("L17").Value2 = 1
$ws.Range("M17").Value2 = 44.68687199999999
$ws.Range("N17").Value2 = 134.060616
$ws.Range("O17").Value2 = 0.3241096696788354
$ws.Range("P17").Value2 = 0.3241096696788355
$ws.Range("Q17").Value2 = 2151.794047805616
$ws.Range("R17").Value2 = 19366.14643025054
$ws.Range("S17").Value2 = 0.1627360114540757
$ws.Range("T17").Value2 = 0.1627360114540757
